$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.039.95'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.423.67'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '411.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.02'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.84%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.90%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '43.46'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.72%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.962.95'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000214'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.81%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.06'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.432.83'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.53'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '61.966.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +20.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '92.20'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.30'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.55'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.36'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '34.35'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.17'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.29%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.77'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.78'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.05'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.168'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.96'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.54'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +10.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0494'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '150.81'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.64%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.137'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.67%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.42'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.326'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.12%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.13'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +7.08%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.96'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.64'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +10.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.25'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.37'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +22.35%  '
$ws.Range('B48').Value = 'Celestia'
$ws.Range('C48').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.59'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.75'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.72%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.148'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +15.53%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '117.40'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +22.17%  '
